# Add a new "05-10-2020" column (T) to the deceased-cases time-series sheet,
# mirroring the existing "04-10-2020" column (S): a styled date header in
# row 1 plus each state/UT's updated cumulative death count for rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell T1 ------------------------------------------------------
# Force a text number format first so the "dd-mm-yyyy"-looking label is
# stored as literal text (matching the other date headers in row 1)
# instead of being auto-converted into a date serial number.
$ws.Range("T1").NumberFormat = "@"
$ws.Range("T1").Value = "05-10-2020"

# Match the visual style used by the rest of the date headers (bold,
# centered, thin box border) the same way S1 ("04-10-2020") is styled.
$ws.Range("T1").Font.Bold = $true
$ws.Range("T1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("T1").VerticalAlignment = -4160     # xlTop
$ws.Range("T1").Borders.LineStyle = 1         # xlContinuous

# --- Per-state cumulative death counts for 05-10-2020 (column T) --------
$values = @{
    2  = 53
    3  = 5981
    4  = 18
    5  = 749
    6  = 915
    7  = 174
    8  = 1045
    9  = 2
    10 = 5510
    11 = 456
    12 = 3496
    13 = 1470
    14 = 217
    15 = 1242
    16 = 743
    17 = 9286
    18 = 836
    19 = 61
    20 = 2434
    21 = 38084
    22 = 74
    23 = 54
    24 = 0
    25 = 17
    26 = 907
    27 = 539
    28 = 3603
    29 = 1545
    30 = 45
    31 = 9784
    32 = 1171
    33 = 299
    34 = 652
    35 = 6029
    36 = 5194
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 20).Value = $values[$row]
}
